$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2330623306233062
$ws.Cells.Item(2, 3).Value = 0.4905149051490515
$ws.Cells.Item(2, 10).Value = 0.008130081300813009
$ws.Cells.Item(2, 16).Value = 0.1598915989159892
$ws.Cells.Item(2, 19).Value = 0.1084010840108401
$ws.Cells.Item(3, 2).Value = 0.005208333333333333
$ws.Cells.Item(3, 3).Value = 0.04166666666666666
$ws.Cells.Item(3, 10).Value = 0.03645833333333334
$ws.Cells.Item(3, 16).Value = 0.7864583333333334
$ws.Cells.Item(3, 19).Value = 0.1302083333333333
$ws.Cells.Item(4, 10).Value = 0.1090909090909091
$ws.Cells.Item(4, 16).Value = 0.4363636363636363
$ws.Cells.Item(4, 19).Value = 0.4545454545454545
$ws.Cells.Item(6, 2).Value = 0.07111111111111111
$ws.Cells.Item(6, 4).Value = 0.008888888888888889
$ws.Cells.Item(6, 6).Value = 0.03555555555555556
$ws.Cells.Item(6, 10).Value = 0.2533333333333334
$ws.Cells.Item(6, 15).Value = 0.008888888888888889
$ws.Cells.Item(6, 17).Value = 0.2088888888888889
$ws.Cells.Item(6, 18).Value = 0.01777777777777778
$ws.Cells.Item(6, 19).Value = 0.3955555555555555
$ws.Cells.Item(7, 2).Value = 0.1106382978723404
$ws.Cells.Item(7, 4).Value = 0.03404255319148936
$ws.Cells.Item(7, 5).Value = 0.00425531914893617
$ws.Cells.Item(7, 6).Value = 0.03829787234042553
$ws.Cells.Item(7, 10).Value = 0.1574468085106383
$ws.Cells.Item(7, 15).Value = 0.02127659574468085
$ws.Cells.Item(7, 17).Value = 0.148936170212766
$ws.Cells.Item(7, 18).Value = 0.06382978723404255
$ws.Cells.Item(7, 19).Value = 0.4212765957446808
$ws.Cells.Item(8, 2).Value = 0.09467455621301775
$ws.Cells.Item(8, 4).Value = 0.01577909270216963
$ws.Cells.Item(8, 5).Value = 0.003944773175542407
$ws.Cells.Item(8, 6).Value = 0.07889546351084813
$ws.Cells.Item(8, 10).Value = 0.1025641025641026
$ws.Cells.Item(8, 15).Value = 0.01972386587771203
$ws.Cells.Item(8, 17).Value = 0.2149901380670611
$ws.Cells.Item(8, 18).Value = 0.04930966469428008
$ws.Cells.Item(8, 19).Value = 0.4201183431952663
$ws.Cells.Item(9, 2).Value = 0.08743169398907104
$ws.Cells.Item(9, 4).Value = 0.02185792349726776
$ws.Cells.Item(9, 6).Value = 0.07103825136612021
$ws.Cells.Item(9, 10).Value = 0.1147540983606557
$ws.Cells.Item(9, 15).Value = 0.02185792349726776
$ws.Cells.Item(9, 17).Value = 0.1748633879781421
$ws.Cells.Item(9, 18).Value = 0.06557377049180328
$ws.Cells.Item(9, 19).Value = 0.4426229508196721
$ws.Cells.Item(10, 2).Value = 0.1220362622036262
$ws.Cells.Item(10, 4).Value = 0.02301255230125523
$ws.Cells.Item(10, 5).Value = 0.002789400278940028
$ws.Cells.Item(10, 6).Value = 0.06066945606694561
$ws.Cells.Item(10, 10).Value = 0.1248256624825663
$ws.Cells.Item(10, 15).Value = 0.01394700139470014
$ws.Cells.Item(10, 17).Value = 0.2412831241283124
$ws.Cells.Item(10, 18).Value = 0.04951185495118549
$ws.Cells.Item(10, 19).Value = 0.3619246861924686
$ws.Cells.Item(11, 7).Value = 0.137466307277628
$ws.Cells.Item(11, 10).Value = 0.09703504043126684
$ws.Cells.Item(11, 11).Value = 0.1752021563342318
$ws.Cells.Item(11, 12).Value = 0.5849056603773585
$ws.Cells.Item(11, 19).Value = 0.005390835579514825
$ws.Cells.Item(12, 7).Value = 0.7431192660550459
$ws.Cells.Item(12, 10).Value = 0.2201834862385321
$ws.Cells.Item(12, 12).Value = 0.009174311926605505
$ws.Cells.Item(12, 19).Value = 0.02752293577981652
$ws.Cells.Item(13, 7).Value = 0.7692307692307693
$ws.Cells.Item(13, 10).Value = 0.1794871794871795
$ws.Cells.Item(13, 19).Value = 0.05128205128205128
$ws.Cells.Item(15, 6).Value = 0.01260504201680672
$ws.Cells.Item(15, 8).Value = 0.1554621848739496
$ws.Cells.Item(15, 9).Value = 0.0546218487394958
$ws.Cells.Item(15, 10).Value = 0.3739495798319328
$ws.Cells.Item(15, 11).Value = 0.08403361344537816
$ws.Cells.Item(15, 13).Value = 0.01680672268907563
$ws.Cells.Item(15, 15).Value = 0.06722689075630252
$ws.Cells.Item(15, 19).Value = 0.2352941176470588
$ws.Cells.Item(16, 6).Value = 0.0170940170940171
$ws.Cells.Item(16, 8).Value = 0.1965811965811966
$ws.Cells.Item(16, 9).Value = 0.07692307692307693
$ws.Cells.Item(16, 10).Value = 0.4017094017094017
$ws.Cells.Item(16, 11).Value = 0.1025641025641026
$ws.Cells.Item(16, 13).Value = 0.01282051282051282
$ws.Cells.Item(16, 15).Value = 0.05982905982905983
$ws.Cells.Item(16, 19).Value = 0.1324786324786325
$ws.Cells.Item(17, 6).Value = 0.01573426573426574
$ws.Cells.Item(17, 8).Value = 0.1660839160839161
$ws.Cells.Item(17, 9).Value = 0.09265734265734266
$ws.Cells.Item(17, 10).Value = 0.4178321678321678
$ws.Cells.Item(17, 11).Value = 0.1241258741258741
$ws.Cells.Item(17, 13).Value = 0.005244755244755245
$ws.Cells.Item(17, 15).Value = 0.06818181818181818
$ws.Cells.Item(17, 19).Value = 0.1101398601398601
$ws.Cells.Item(18, 6).Value = 0.03149606299212598
$ws.Cells.Item(18, 8).Value = 0.1102362204724409
$ws.Cells.Item(18, 9).Value = 0.1023622047244094
$ws.Cells.Item(18, 10).Value = 0.4330708661417323
$ws.Cells.Item(18, 11).Value = 0.1417322834645669
$ws.Cells.Item(18, 13).Value = 0.02362204724409449
$ws.Cells.Item(18, 15).Value = 0.07874015748031496
$ws.Cells.Item(18, 19).Value = 0.07874015748031496
$ws.Cells.Item(19, 6).Value = 0.0104602510460251
$ws.Cells.Item(19, 8).Value = 0.2231520223152022
$ws.Cells.Item(19, 9).Value = 0.06276150627615062
$ws.Cells.Item(19, 10).Value = 0.3640167364016736
$ws.Cells.Item(19, 11).Value = 0.1164574616457462
$ws.Cells.Item(19, 13).Value = 0.02092050209205021
$ws.Cells.Item(19, 15).Value = 0.06136680613668061
$ws.Cells.Item(19, 19).Value = 0.1408647140864714
